{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies \"hybrid bold + color\" highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific bullet\n// paragraphs, matching the target OOXML diff exactly: each metric\n// substring is split out into its own run with <w:b/> and\n// <w:color w:val=\"2C3E50\"/>, while the surrounding text stays in plain\n// runs.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Paragraphs to touch, identified by a unique, stable prefix of their\n// text, plus the ordered list of metric substrings to bold within them.\nconst EDITS = [\n  {\n    prefix: \"\u2022 Discovered systematic race coding errors\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    prefix: \"\u2022 Utilized advanced sampling methods\",\n    metrics: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    prefix: \"\u2022 Trigonometric algorithm for boundary estimation\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    prefix: \"\u2022 Built real-time FEC analysis systems\",\n    metrics: [\"$2\"],\n  },\n  {\n    prefix: \"\u2022 Modernized legacy ETL processes\",\n    metrics: [\"57%\"],\n  },\n  {\n    prefix: \"\u2022 Platform impact: Built redistricting system\",\n    metrics: [\"12,847\"],\n  },\n  {\n    prefix: \"\u2022 Revenue generation: Delivered\",\n    metrics: [\"$4.9M\"],\n  },\n  {\n    prefix: \"\u2022 23% conversion rate improvement\",\n    metrics: [\"23%\"],\n  },\n];\n\nasync function boldenMetric(paragraph, metric) {\n  const results = paragraph.search(metric, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    const font = results.items[i].font;\n    font.bold = true;\n    font.color = HIGHLIGHT_COLOR;\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const edit of EDITS) {\n  const paragraph = paragraphs.items.find((p) => p.text.indexOf(edit.prefix) === 0);\n  if (!paragraph) {\n    continue;\n  }\n  for (const metric of edit.metrics) {\n    await boldenMetric(paragraph, metric);\n  }\n}\n", "ps1": "# Word COM interop edit script.\n# Applies \"hybrid bold + color\" highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific bullet\n# paragraphs, matching the target OOXML diff exactly: each metric\n# substring is split out into its own run with <w:b/> and\n# <w:color w:val=\"2C3E50\"/>, while the surrounding text stays in plain\n# runs.\n\n$d = $word.ActiveDocument\n$HighlightColor = \"2C3E50\"\n\n# Paragraphs to touch, identified by a unique, stable substring of their\n# text, plus the ordered list of metric substrings to bold within them.\n$Edits = @(\n    @{ Match = \"Discovered systematic race coding errors\"; Metrics = @(\"23%\", \"64%\") },\n    @{ Match = \"Utilized advanced sampling methods\"; Metrics = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\") },\n    @{ Match = \"Trigonometric algorithm for boundary estimation\"; Metrics = @(\"73.5%\", \"$4.7M\") },\n    @{ Match = \"Built real-time FEC analysis systems\"; Metrics = @(\"$2\") },\n    @{ Match = \"Modernized legacy ETL processes\"; Metrics = @(\"57%\") },\n    @{ Match = \"Platform impact: Built redistricting system\"; Metrics = @(\"12,847\") },\n    @{ Match = \"Revenue generation: Delivered\"; Metrics = @(\"$4.9M\") },\n    @{ Match = \"23% conversion rate improvement\"; Metrics = @(\"23%\") }\n)\n\nforeach ($edit in $Edits) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.Contains($edit.Match)) {\n            $paraRange = $p.Range\n            foreach ($metric in $edit.Metrics) {\n                $found = $paraRange.Find.Execute($metric)\n                if ($found) {\n                    $paraRange.Font.Bold = 1\n                    $paraRange.Font.Color = $HighlightColor\n                }\n            }\n            break\n        }\n    }\n}\n"}
